# Rename the inline picture "docPr" names so they match the IDs used by the
# other half of the logo pair:
#   - footer1.xml  (first-page footer,  wp:docPr id="3") image1.png -> image2.png
#   - footer2.xml  (default footer,     wp:docPr id="2") image1.png -> image2.png
#   - header1.xml  (first-page header,  wp:docPr id="1") image2.jpg -> image1.jpg
#
# WdHeaderFooterIndex: 1 = wdHeaderFooterPrimary (default), 2 = wdHeaderFooterFirstPage
#   Footers.Item(1) -> footer2.xml (default footer)    -> image1.png -> image2.png
#   Footers.Item(2) -> footer1.xml (first-page footer) -> image1.png -> image2.png
#   Headers.Item(2) -> header1.xml (first-page header) -> image2.jpg -> image1.jpg

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

$defaultFooterLogo = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$defaultFooterLogo.Name = "image2.png"
Write-Output "Renamed default footer logo to $($defaultFooterLogo.Name)"

$firstPageFooterLogo = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$firstPageFooterLogo.Name = "image2.png"
Write-Output "Renamed first-page footer logo to $($firstPageFooterLogo.Name)"

$firstPageHeaderLogo = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$firstPageHeaderLogo.Name = "image1.jpg"
Write-Output "Renamed first-page header logo to $($firstPageHeaderLogo.Name)"
